$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header labels: reorder block names
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "living_rooms_1"
$ws.Range("F1").Value = "living_rooms_2"

# Row 2: move the "1" marker from E2 to C2
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = 0

# Row 4: move the "1" marker from C4 to F4
$ws.Range("C4").Value = 0
$ws.Range("F4").Value = 1

# Row 5: move the "1" marker from D5 to E5
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1

# Row 6: move the "1" marker from F6 to B6
$ws.Range("B6").Value = 1
$ws.Range("F6").Value = 0

# Row 7: move the "1" marker from B7 to D7
$ws.Range("B7").Value = 0
$ws.Range("D7").Value = 1
